$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2446
$ws.Range("I62").Value = 2446
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2446
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1822
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 2446
$ws.Range("I65").Value = 2446
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 12230
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -9110
$ws.Range("N65").ClearContents()

$ws.Range("H108").Value = 60000
$ws.Range("J108").Value = 60000
$ws.Range("L108").Value = 60000
$ws.Range("N108").Value = -67680

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H112").Value = 1955.375
$ws.Range("J112").Value = 2002.3226
$ws.Range("L112").Value = 6006.9678
$ws.Range("N112").Value = -8222.9678

$ws.Range("H116").Value = 12387.5
$ws.Range("J116").Value = 4465
$ws.Range("L116").Value = 4465
$ws.Range("N116").Value = -11349

$ws.Range("H132").Value = 1330.8077
$ws.Range("I132").Value = 1322.5294
$ws.Range("K132").Value = 3967.5882
$ws.Range("M132").Value = -1437.5882

$ws.Range("H138").Value = 2079.2415
$ws.Range("I138").Value = 1840.1082
$ws.Range("J138").Value = 2500.5715
$ws.Range("K138").Value = 5520.3246
$ws.Range("L138").Value = 7501.7145
$ws.Range("M138").Value = -380.3245999999999
$ws.Range("N138").Value = -17781.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5156.256
$ws.Range("I32").Value = 4006.3057
$ws.Range("K32").Value = 4006.3057
$ws.Range("M32").Value = -3719.3057

$ws.Range("H74").Value = 1094.2646
$ws.Range("I74").Value = 586.1429000000001
$ws.Range("J74").Value = 3465.5
$ws.Range("K74").Value = 586.1429000000001
$ws.Range("L74").Value = 3465.5
$ws.Range("M74").Value = 287.8570999999999
$ws.Range("N74").Value = -5213.5

$ws.Range("H77").Value = 1094.2646
$ws.Range("I77").Value = 586.1429000000001
$ws.Range("J77").Value = 3465.5
$ws.Range("K77").Value = 2930.7145
$ws.Range("L77").Value = 17327.5
$ws.Range("M77").Value = 1437.2855
$ws.Range("N77").Value = -26063.5

$ws.Range("H132").Value = 1779.1111
$ws.Range("I132").Value = 1712.1428
$ws.Range("K132").Value = 5136.428400000001
$ws.Range("M132").Value = -2606.428400000001

$ws.Range("H139").Value = 45000
$ws.Range("J139").Value = 45000
$ws.Range("L139").Value = 45000
$ws.Range("N139").Value = -55280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3270.2222
$ws.Range("I20").Value = 3286.4
$ws.Range("K20").Value = 3286.4
$ws.Range("M20").Value = -3039.4

$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()

$ws.Range("H86").Value = 860094.9399999999
$ws.Range("I86").Value = 1003739.8
$ws.Range("J86").Value = 500982.75
$ws.Range("K86").Value = 1003739.8
$ws.Range("L86").Value = 500982.75
$ws.Range("M86").Value = -1002616.8
$ws.Range("N86").Value = -503228.75

$ws.Range("H89").Value = 860094.9399999999
$ws.Range("I89").Value = 1003739.8
$ws.Range("J89").Value = 500982.75
$ws.Range("K89").Value = 5018699
$ws.Range("L89").Value = 2504913.75
$ws.Range("M89").Value = -5013083
$ws.Range("N89").Value = -2516145.75

$ws.Range("H94").Value = 872.5833
$ws.Range("I94").Value = 855.44446
$ws.Range("K94").Value = 855.44446
$ws.Range("M94").Value = -404.44446

$ws.Range("H134").Value = 13113.637
$ws.Range("J134").Value = 2354.75
$ws.Range("L134").Value = 7064.25
$ws.Range("N134").Value = -12134.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 67445.55499999999
$ws.Range("I16").Value = 86487.14
$ws.Range("K16").Value = 86487.14
$ws.Range("M16").Value = -86200.14

$ws.Range("H113").Value = 67445.55499999999
$ws.Range("I113").Value = 86487.14
$ws.Range("K113").Value = 86487.14
$ws.Range("M113").Value = -84317.14

$ws.Range("H122").Value = 3601.3333
$ws.Range("I122").Value = 2818.9092
$ws.Range("J122").Value = 5753
$ws.Range("K122").Value = 8456.7276
$ws.Range("L122").Value = 17259
$ws.Range("M122").Value = -6006.7276
$ws.Range("N122").Value = -22159

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 5822
$ws.Range("J106").Value = 5822
$ws.Range("L106").Value = 17466
$ws.Range("N106").Value = -19358

$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("M108").ClearContents()
$ws.Range("N108").ClearContents()

$ws.Range("H131").Value = 30320.125
$ws.Range("J131").Value = 33001.5
$ws.Range("L131").Value = 99004.5
$ws.Range("N131").Value = -109084.5

$ws.Range("H132").Value = 1780.8
$ws.Range("J132").Value = 2333
$ws.Range("L132").Value = 20997
$ws.Range("N132").Value = -26057

$ws.Range("H140").Value = 2459.6316
$ws.Range("I140").Value = 1549.091
$ws.Range("J140").Value = 3711.625
$ws.Range("K140").Value = 4647.272999999999
$ws.Range("L140").Value = 11134.875
$ws.Range("M140").Value = 532.7270000000008
$ws.Range("N140").Value = -21494.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 999
$ws.Range("I80").Value = 999
$ws.Range("K80").Value = 999
$ws.Range("M80").Value = -1

$ws.Range("H83").Value = 999
$ws.Range("I83").Value = 999
$ws.Range("K83").Value = 4995
$ws.Range("M83").Value = -3

$ws.Range("H122").Value = 1902.3334
$ws.Range("I122").Value = 1538
$ws.Range("K122").Value = 4614
$ws.Range("M122").Value = -2164

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3280.5557
$ws.Range("I7").Value = 3230
$ws.Range("J7").Value = 3343.75
$ws.Range("K7").Value = 3230
$ws.Range("L7").Value = 3343.75
$ws.Range("M7").Value = -3118
$ws.Range("N7").Value = -3567.75

$ws.Range("H40").Value = 13214.846
$ws.Range("I40").Value = 17428.857
$ws.Range("J40").Value = 8298.5
$ws.Range("K40").Value = 17428.857
$ws.Range("L40").Value = 8298.5
$ws.Range("M40").Value = -17292.857
$ws.Range("N40").Value = -8570.5

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H126").Value = 3280.5557
$ws.Range("I126").Value = 3230
$ws.Range("J126").Value = 3343.75
$ws.Range("K126").Value = 9690
$ws.Range("L126").Value = 10031.25
$ws.Range("M126").Value = -7220
$ws.Range("N126").Value = -14971.25

$ws.Range("H132").Value = 2456.35
$ws.Range("I132").Value = 1550.5
$ws.Range("J132").Value = 2844.5715
$ws.Range("K132").Value = 4651.5
$ws.Range("L132").Value = 8533.7145
$ws.Range("M132").Value = -2121.5
$ws.Range("N132").Value = -13593.7145

$ws.Range("H136").Value = 3659.923
$ws.Range("I136").Value = 2924.4783
$ws.Range("K136").Value = 8773.4349
$ws.Range("M136").Value = -6223.4349

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 38999.668
$ws.Range("J108").Value = 38999.668
$ws.Range("L108").Value = 38999.668
$ws.Range("N108").Value = -46679.668

$ws.Range("H122").Value = 35035.086
$ws.Range("I122").Value = 56677
$ws.Range("K122").Value = 170031
$ws.Range("M122").Value = -167581

$ws.Range("H126").Value = 8486.0625
$ws.Range("I126").Value = 8898.429
$ws.Range("J126").Value = 5599.5
$ws.Range("K126").Value = 26695.287
$ws.Range("L126").Value = 16798.5
$ws.Range("M126").Value = -24225.287
$ws.Range("N126").Value = -21738.5

$ws.Range("H132").Value = 2010.6522
$ws.Range("I132").Value = 1316.0588
$ws.Range("K132").Value = 3948.1764
$ws.Range("M132").Value = -1418.1764

$ws.Range("H136").Value = 13229258
$ws.Range("I136").Value = 21369194
$ws.Range("J136").Value = 1863.0625
$ws.Range("K136").Value = 64107582
$ws.Range("L136").Value = 5589.1875
$ws.Range("M136").Value = -64105032
$ws.Range("N136").Value = -10689.1875

